# chore: update Sheets via scheduled runner
# Refresh computed market-price / leve-profit figures across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 90.2
$ws.Range("I2").Value = 87.75
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 87.75
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 25.25
$ws.Range("N2").Value = -326
$ws.Range("H74").Value = 5212269.5
$ws.Range("I74").Value = 3300.6
$ws.Range("J74").Value = 6583051
$ws.Range("K74").Value = 3300.6
$ws.Range("L74").Value = 6583051
$ws.Range("M74").Value = -2364.6
$ws.Range("N74").Value = -6584923
$ws.Range("H77").Value = 5212269.5
$ws.Range("I77").Value = 3300.6
$ws.Range("J77").Value = 6583051
$ws.Range("K77").Value = 16503
$ws.Range("L77").Value = 32915255
$ws.Range("M77").Value = -11823
$ws.Range("N77").Value = -32924615
$ws.Range("H98").Value = 1074.25
$ws.Range("I98").Value = 999.6667
$ws.Range("K98").Value = 999.6667
$ws.Range("M98").Value = 498.3333
$ws.Range("H106").Value = 1095.409
$ws.Range("I106").Value = 814.45
$ws.Range("K106").Value = 814.45
$ws.Range("M106").Value = -183.45
$ws.Range("H107").Value = 1657.1428
$ws.Range("I107").Value = 1920
$ws.Range("K107").Value = 1920
$ws.Range("M107").Value = 0
$ws.Range("H116").Value = 4258.8125
$ws.Range("I116").Value = 1792.7778
$ws.Range("K116").Value = 1792.7778
$ws.Range("M116").Value = 1649.2222
$ws.Range("H122").Value = 1074.25
$ws.Range("I122").Value = 999.6667
$ws.Range("K122").Value = 2999.0001
$ws.Range("M122").Value = -549.0001000000002
$ws.Range("H129").Value = 858.4666999999999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 858.4666999999999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 2575.4001
$ws.Range("M129").Value = $null
$ws.Range("N129").Value = -12575.4001
$ws.Range("H138").Value = 1757.7291
$ws.Range("I138").Value = 546.6786
$ws.Range("J138").Value = 3453.2
$ws.Range("K138").Value = 1640.0358
$ws.Range("L138").Value = 10359.6
$ws.Range("M138").Value = 3499.9642
$ws.Range("N138").Value = -20639.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2026.65
$ws.Range("I2").Value = 1396.5
$ws.Range("J2").Value = 3497
$ws.Range("K2").Value = 1396.5
$ws.Range("L2").Value = 3497
$ws.Range("M2").Value = -1283.5
$ws.Range("N2").Value = -3723
$ws.Range("H32").Value = 26180.861
$ws.Range("I32").Value = 27919.924
$ws.Range("K32").Value = 27919.924
$ws.Range("M32").Value = -27632.924
$ws.Range("H74").Value = 2176.1853
$ws.Range("I74").Value = 2150.28
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 2150.28
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -1276.28
$ws.Range("N74").Value = -4248
$ws.Range("H77").Value = 2176.1853
$ws.Range("I77").Value = 2150.28
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 10751.4
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -6383.400000000001
$ws.Range("N77").Value = -21236
$ws.Range("H110").Value = 300
$ws.Range("I110").Value = 300
$ws.Range("K110").Value = 300
$ws.Range("M110").Value = 1745
$ws.Range("H116").Value = 2026.65
$ws.Range("I116").Value = 1396.5
$ws.Range("J116").Value = 3497
$ws.Range("K116").Value = 1396.5
$ws.Range("L116").Value = 3497
$ws.Range("M116").Value = 897.5
$ws.Range("N116").Value = -8085
$ws.Range("H122").Value = 2051.375
$ws.Range("I122").Value = 2130.1428
$ws.Range("K122").Value = 6390.428400000001
$ws.Range("M122").Value = -3940.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2026.65
$ws.Range("I3").Value = 1396.5
$ws.Range("J3").Value = 3497
$ws.Range("K3").Value = 1396.5
$ws.Range("L3").Value = 3497
$ws.Range("M3").Value = -1282.5
$ws.Range("N3").Value = -3725
$ws.Range("H76").Value = 30000
$ws.Range("J76").Value = 30000
$ws.Range("L76").Value = 30000
$ws.Range("N76").Value = -30630
$ws.Range("H79").Value = 30000
$ws.Range("J79").Value = 30000
$ws.Range("L79").Value = 30000
$ws.Range("N79").Value = -32184
$ws.Range("H86").Value = 2088.5557
$ws.Range("I86").Value = 1856.7142
$ws.Range("J86").Value = 2900
$ws.Range("K86").Value = 1856.7142
$ws.Range("L86").Value = 2900
$ws.Range("M86").Value = -733.7141999999999
$ws.Range("N86").Value = -5146
$ws.Range("H89").Value = 2088.5557
$ws.Range("I89").Value = 1856.7142
$ws.Range("J89").Value = 2900
$ws.Range("K89").Value = 9283.571
$ws.Range("L89").Value = 14500
$ws.Range("M89").Value = -3667.571
$ws.Range("N89").Value = -25732
$ws.Range("H94").Value = 2589.5676
$ws.Range("I94").Value = 1329.1923
$ws.Range("J94").Value = 5568.636
$ws.Range("K94").Value = 1329.1923
$ws.Range("L94").Value = 5568.636
$ws.Range("M94").Value = -878.1922999999999
$ws.Range("N94").Value = -6470.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14960470
$ws.Range("I99").Value = 2978594
$ws.Range("K99").Value = 2978594
$ws.Range("M99").Value = -2977096
$ws.Range("H122").Value = 834.0833
$ws.Range("I122").Value = 834.0833
$ws.Range("K122").Value = 2502.2499
$ws.Range("M122").Value = -52.2498999999998
$ws.Range("H126").Value = 14960470
$ws.Range("I126").Value = 2978594
$ws.Range("K126").Value = 8935782
$ws.Range("M126").Value = -8933312

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 766.28284
$ws.Range("J131").Value = 778.3542
$ws.Range("L131").Value = 2335.0626
$ws.Range("N131").Value = -12415.0626

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 425
$ws.Range("J22").Value = 425
$ws.Range("L22").Value = 425
$ws.Range("N22").Value = -1483
$ws.Range("H102").Value = 1739.32
$ws.Range("I102").Value = 1781.9546
$ws.Range("J102").Value = 1426.6666
$ws.Range("K102").Value = 1781.9546
$ws.Range("L102").Value = 1426.6666
$ws.Range("M102").Value = -159.9546
$ws.Range("N102").Value = -4670.6666
$ws.Range("H113").Value = 3531.7144
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 3744.4
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3744.4
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8084.4
$ws.Range("H122").Value = 2066.5
$ws.Range("I122").Value = 2136.8
$ws.Range("J122").Value = 1996.2
$ws.Range("K122").Value = 6410.400000000001
$ws.Range("L122").Value = 5988.6
$ws.Range("M122").Value = -3960.400000000001
$ws.Range("N122").Value = -10888.6
$ws.Range("H126").Value = 4944.4443
$ws.Range("I126").Value = 3357.1428
$ws.Range("J126").Value = 6653.846
$ws.Range("K126").Value = 10071.4284
$ws.Range("L126").Value = 19961.538
$ws.Range("M126").Value = -7601.428400000001
$ws.Range("N126").Value = -24901.538
$ws.Range("H132").Value = 219714.14
$ws.Range("I132").Value = 340333.66
$ws.Range("J132").Value = 129249.5
$ws.Range("K132").Value = 1021000.98
$ws.Range("L132").Value = 387748.5
$ws.Range("M132").Value = -1018470.98
$ws.Range("N132").Value = -392808.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3050.75
$ws.Range("J22").Value = 734
$ws.Range("L22").Value = 734
$ws.Range("N22").Value = -1324
$ws.Range("H27").Value = 3050.75
$ws.Range("J27").Value = 734
$ws.Range("L27").Value = 734
$ws.Range("N27").Value = -948
$ws.Range("H46").Value = 2585.8572
$ws.Range("I46").Value = 2550
$ws.Range("J46").Value = 2600.2
$ws.Range("K46").Value = 2550
$ws.Range("L46").Value = 2600.2
$ws.Range("M46").Value = -2362
$ws.Range("N46").Value = -2976.2
$ws.Range("H93").Value = 1961.1052
$ws.Range("I93").Value = 1909.4706
$ws.Range("J93").Value = 2400
$ws.Range("K93").Value = 1909.4706
$ws.Range("L93").Value = 2400
$ws.Range("M93").Value = -661.4706000000001
$ws.Range("N93").Value = -4896
$ws.Range("H122").Value = 1156916.9
$ws.Range("I122").Value = 2453568
$ws.Range("K122").Value = 7360704
$ws.Range("M122").Value = -7358254

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 12000
$ws.Range("J75").Value = 12000
$ws.Range("L75").Value = 12000
$ws.Range("N75").Value = -13872
$ws.Range("H78").Value = 12000
$ws.Range("J78").Value = 12000
$ws.Range("L78").Value = 36000
$ws.Range("N78").Value = -45360
$ws.Range("H107").Value = 2526327
$ws.Range("I107").Value = 524.1
$ws.Range("J107").Value = 5683580.5
$ws.Range("K107").Value = 1572.3
$ws.Range("L107").Value = 17050741.5
$ws.Range("M107").Value = 347.6999999999998
$ws.Range("N107").Value = -17054581.5
$ws.Range("H122").Value = 1946.9445
$ws.Range("I122").Value = 1797
$ws.Range("J122").Value = 2471.75
$ws.Range("K122").Value = 5391
$ws.Range("L122").Value = 7415.25
$ws.Range("M122").Value = -2941
$ws.Range("N122").Value = -12315.25
$ws.Range("H126").Value = 2075
$ws.Range("I126").Value = 867
$ws.Range("J126").Value = 2799.8
$ws.Range("K126").Value = 2601
$ws.Range("L126").Value = 8399.400000000001
$ws.Range("M126").Value = -131
$ws.Range("N126").Value = -13339.4
